# Update odds values in Sheet1 to match the target revision of
# Jogos_da_Semana_FlashScore_2024-10-16.xlsx (odds refreshed for several matches).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53

# Row 3
$ws.Range("G3").Value = 1.71
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 5.5
$ws.Range("N3").Value = 8.5
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 8.5
$ws.Range("AF3").Value = 67
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 15
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 9.5
$ws.Range("AQ3").Value = 34
$ws.Range("AU3").Value = 9
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 29
$ws.Range("AZ3").Value = 101
$ws.Range("BB3").Value = 351

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("J4").Value = 4.33
$ws.Range("K4").Value = 1.92
$ws.Range("L4").Value = 3
$ws.Range("AC4").Value = 6.5
$ws.Range("AO4").Value = 21
$ws.Range("AX4").Value = 13
$ws.Range("BA4").Value = 81

# Row 5
$ws.Range("G5").Value = 1.51
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48

# Row 6
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("Q6").Value = 2.4
$ws.Range("R6").Value = 1.53

# Row 8
$ws.Range("K8").Value = 2.25
$ws.Range("R8").Value = 2.12
$ws.Range("U8").Value = 1.47
$ws.Range("V8").Value = 2.32
$ws.Range("Y8").Value = 9.5
$ws.Range("AD8").Value = 7.2
$ws.Range("AH8").Value = 11
$ws.Range("AI8").Value = 14.5
$ws.Range("AM8").Value = 22
$ws.Range("AO8").Value = 12.5
$ws.Range("AR8").Value = 65
$ws.Range("AT8").Value = 2.9
$ws.Range("AU8").Value = 6.3
$ws.Range("AY8").Value = 17
$ws.Range("BA8").Value = 65

# Row 11
$ws.Range("K11").Value = 1.87
$ws.Range("M11").Value = 1.11
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 2.38
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.53

# Row 12
$ws.Range("G12").Value = 1.57
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 6.5
$ws.Range("J12").Value = 2.3
$ws.Range("K12").Value = 2.05
$ws.Range("L12").Value = 6.5
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.57
$ws.Range("S12").Value = 1.5
$ws.Range("T12").Value = 2.5
$ws.Range("U12").Value = 2.38
$ws.Range("V12").Value = 1.53
$ws.Range("Y12").Value = 9
$ws.Range("Z12").Value = 11
$ws.Range("AC12").Value = 7
$ws.Range("AD12").Value = 7
$ws.Range("AE12").Value = 23
$ws.Range("AK12").Value = 81
$ws.Range("AO12").Value = 8.5
$ws.Range("AQ12").Value = 29
$ws.Range("AR12").Value = 51
$ws.Range("AT12").Value = 2.5
$ws.Range("AW12").Value = 7.5
$ws.Range("AZ12").Value = 151

# Row 13
$ws.Range("H13").Value = 3.25
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("O13").Value = 1.29
$ws.Range("Q13").Value = 1.92
$ws.Range("R13").Value = 1.82
$ws.Range("U13").Value = 1.73
$ws.Range("Y13").Value = 13
$ws.Range("AC13").Value = 10
$ws.Range("AQ13").Value = 67
$ws.Range("AS13").Value = 201

# Row 14
$ws.Range("G14").Value = 2.32
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.8
$ws.Range("J14").Value = 2.95
$ws.Range("L14").Value = 3.45
$ws.Range("S14").Value = 1.44
$ws.Range("T14").Value = 2.62
$ws.Range("W14").Value = 7.3
$ws.Range("X14").Value = 11
$ws.Range("Y14").Value = 9.25
$ws.Range("Z14").Value = 24
$ws.Range("AA14").Value = 20
$ws.Range("AH14").Value = 8
$ws.Range("AI14").Value = 13.5
$ws.Range("AJ14").Value = 10.75
$ws.Range("AK14").Value = 35
$ws.Range("AL14").Value = 26
$ws.Range("AM14").Value = 37
$ws.Range("AN14").Value = 4.2
$ws.Range("AO14").Value = 12.5
$ws.Range("AP14").Value = 22
$ws.Range("AQ14").Value = 50
$ws.Range("AR14").Value = 90
$ws.Range("AT14").Value = 2.62
$ws.Range("AW14").Value = 4.7
$ws.Range("AX14").Value = 15.5
$ws.Range("AY14").Value = 25
$ws.Range("AZ14").Value = 75
